# Automatic update of files.
# The underlying records got re-keyed: the record formerly shown on row 2
# moved to row 4, the record formerly on row 3 moved to row 2, and the
# record formerly on row 4 moved to row 3 (a 3-way cyclic rotation of the
# data rows). Row numbers/labels themselves stay put; only the field
# values move between rows. Apply that rotation explicitly per cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- snapshot the "before" values we need, from their original rows ----

# Row 2 (was: Korallrot / Corallorhiza trifida / Anders Granér ...)
$A2 = $ws.Range("A2").Value()
$B2 = $ws.Range("B2").Value()
$E2 = $ws.Range("E2").Value()
$F2 = $ws.Range("F2").Value()
$G2 = $ws.Range("G2").Value()
$H2 = $ws.Range("H2").Value()
$P2 = $ws.Range("P2").Value()
$Q2 = $ws.Range("Q2").Value()
$R2 = $ws.Range("R2").Value()
$AX2 = $ws.Range("AX2").Value()

# Row 3 (was: Fläcknycklar / Dactylorhiza maculata, id 106541918 ...)
$A3 = $ws.Range("A3").Value()
$B3 = $ws.Range("B3").Value()
$E3 = $ws.Range("E3").Value()
$F3 = $ws.Range("F3").Value()
$G3 = $ws.Range("G3").Value()
$H3 = $ws.Range("H3").Value()
$P3 = $ws.Range("P3").Value()
$Q3 = $ws.Range("Q3").Value()
$R3 = $ws.Range("R3").Value()
$AX3 = $ws.Range("AX3").Value()

# Row 4 (was: Fläcknycklar / Dactylorhiza maculata, id 106541919 ...)
$A4 = $ws.Range("A4").Value()
$B4 = $ws.Range("B4").Value()
$E4 = $ws.Range("E4").Value()
$F4 = $ws.Range("F4").Value()
$G4 = $ws.Range("G4").Value()
$H4 = $ws.Range("H4").Value()
$P4 = $ws.Range("P4").Value()
$Q4 = $ws.Range("Q4").Value()
$R4 = $ws.Range("R4").Value()
$AX4 = $ws.Range("AX4").Value()

# ---- write the rotated values back: 2 <- 3, 3 <- 4, 4 <- 2 ----

# New row 2 gets old row 3's full record.
$ws.Range("A2").Value = $A3
$ws.Range("B2").Value = $B3
$ws.Range("E2").Value = $E3
$ws.Range("F2").Value = $F3
$ws.Range("G2").Value = $G3
$ws.Range("H2").Value = $H3
$ws.Range("P2").Value = $P3
$ws.Range("Q2").Value = $Q3
$ws.Range("R2").Value = $R3
$ws.Range("AX2").Value = $AX3

# New row 3 gets old row 4's full record.
$ws.Range("A3").Value = $A4
$ws.Range("B3").Value = $B4
$ws.Range("E3").Value = $E4
$ws.Range("F3").Value = $F4
$ws.Range("G3").Value = $G4
$ws.Range("H3").Value = $H4
$ws.Range("P3").Value = $P4
$ws.Range("Q3").Value = $Q4
$ws.Range("R3").Value = $R4
$ws.Range("AX3").Value = $AX4

# New row 4 gets old row 2's full record.
$ws.Range("A4").Value = $A2
$ws.Range("B4").Value = $B2
$ws.Range("E4").Value = $E2
$ws.Range("F4").Value = $F2
$ws.Range("G4").Value = $G2
$ws.Range("H4").Value = $H2
$ws.Range("P4").Value = $P2
$ws.Range("Q4").Value = $Q2
$ws.Range("R4").Value = $R2
$ws.Range("AX4").Value = $AX2
